# Prefix each worksheet's protocol/sheet name onto the "Step/label" values
# found in column A (rows 2..N), for the set of protocol sheets that were
# previously sharing duplicate, non-unique command names across tabs.
#
# Example: on sheet "free1", cell A2 "Step4 Seed" becomes "free1 Step4 Seed".

$wb = $excel.ActiveWorkbook

$targetSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $firstRow = $usedRange.Row
    $lastRow = $firstRow + $usedRange.Rows.Count - 1

    # Row 1 is the header ("Name"/"Text"/"Note"/"*Guidelines"); data starts row 2.
    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, 1)
        $current = $cell.Text

        if ([string]::IsNullOrEmpty($current)) {
            continue
        }

        $prefix = $sheetName + " "
        if ($current.StartsWith($prefix)) {
            continue
        }

        $cell.Value = $prefix + $current
    }
}
